$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 540.7273
$ws.Range("I11").Value = 540.7273
$ws.Range("K11").Value = 540.7273
$ws.Range("M11").Value = -400.7273
$ws.Range("H19").Value = 1489.5
$ws.Range("I19").Value = 1489
$ws.Range("K19").Value = 1489
$ws.Range("M19").Value = -1314
$ws.Range("H33").Value = 398.82352
$ws.Range("I33").Value = 184.61539
$ws.Range("K33").Value = 184.61539
$ws.Range("M33").Value = 44.38461000000001
$ws.Range("H70").Value = 5988
$ws.Range("I70").Value = 2073
$ws.Range("J70").Value = 7106.5713
$ws.Range("K70").Value = 6219
$ws.Range("L70").Value = 21319.7139
$ws.Range("M70").Value = -5949
$ws.Range("N70").Value = -21859.7139
$ws.Range("H73").Value = 5988
$ws.Range("I73").Value = 2073
$ws.Range("J73").Value = 7106.5713
$ws.Range("K73").Value = 6219
$ws.Range("L73").Value = 21319.7139
$ws.Range("M73").Value = -5283
$ws.Range("N73").Value = -23191.7139
$ws.Range("H125").Value = 2780.4614
$ws.Range("J125").Value = 2893.375
$ws.Range("L125").Value = 26040.375
$ws.Range("N125").Value = -30960.375
$ws.Range("H129").Value = 1121.6875
$ws.Range("I129").Value = 837.125
$ws.Range("K129").Value = 2511.375
$ws.Range("M129").Value = 2488.625
$ws.Range("H132").Value = 5650.4443
$ws.Range("I132").Value = 4820.793
$ws.Range("K132").Value = 14462.379
$ws.Range("M132").Value = -11932.379
$ws.Range("H133").Value = 84622.60000000001
$ws.Range("J133").Value = 84622.60000000001
$ws.Range("L133").Value = 84622.60000000001
$ws.Range("N133").Value = -94742.60000000001
$ws.Range("H135").Value = 1075.2667
$ws.Range("I135").Value = 1009.4
$ws.Range("K135").Value = 9084.6
$ws.Range("M135").Value = -6549.6
$ws.Range("H137").Value = 1581.091
$ws.Range("I137").Value = 1052
$ws.Range("J137").Value = 2022
$ws.Range("K137").Value = 3156
$ws.Range("L137").Value = 6066
$ws.Range("M137").Value = -606
$ws.Range("N137").Value = -11166
$ws.Range("H138").Value = 2491.2903
$ws.Range("I138").Value = 2098.6
$ws.Range("J138").Value = 2678.2856
$ws.Range("K138").Value = 6295.799999999999
$ws.Range("L138").Value = 8034.8568
$ws.Range("M138").Value = -1155.799999999999
$ws.Range("N138").Value = -18314.8568
$ws.Range("H141").Value = 4248.9443
$ws.Range("I141").Value = 2460.4666
$ws.Range("K141").Value = 7381.399800000001
$ws.Range("M141").Value = -2201.399800000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 500
$ws.Range("I3").Value = 500
$ws.Range("K3").Value = 500
$ws.Range("M3").Value = -385
$ws.Range("H5").Value = 430
$ws.Range("I5").Value = 406
$ws.Range("J5").Value = 450
$ws.Range("K5").Value = 406
$ws.Range("L5").Value = 450
$ws.Range("M5").Value = -294
$ws.Range("N5").Value = -674
$ws.Range("H32").Value = 18688.541
$ws.Range("I32").Value = 24428.418
$ws.Range("K32").Value = 24428.418
$ws.Range("M32").Value = -24141.418
$ws.Range("H74").Value = 2820.5386
$ws.Range("I74").Value = 2913.818
$ws.Range("J74").Value = 2307.5
$ws.Range("K74").Value = 2913.818
$ws.Range("L74").Value = 2307.5
$ws.Range("M74").Value = -2039.818
$ws.Range("N74").Value = -4055.5
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H77").Value = 2820.5386
$ws.Range("I77").Value = 2913.818
$ws.Range("J77").Value = 2307.5
$ws.Range("K77").Value = 14569.09
$ws.Range("L77").Value = 11537.5
$ws.Range("M77").Value = -10201.09
$ws.Range("N77").Value = -20273.5
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H122").Value = 2235.9
$ws.Range("I122").Value = 1889.3529
$ws.Range("K122").Value = 5668.0587
$ws.Range("M122").Value = -3218.0587
$ws.Range("H132").Value = 35743.363
$ws.Range("I132").Value = 94155.17999999999
$ws.Range("J132").Value = 6537.4546
$ws.Range("K132").Value = 282465.54
$ws.Range("L132").Value = 19612.3638
$ws.Range("M132").Value = -279935.54
$ws.Range("N132").Value = -24672.3638
$ws.Range("H135").Value = 52500
$ws.Range("J135").Value = 52500
$ws.Range("L135").Value = 52500
$ws.Range("N135").Value = -62640

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 430
$ws.Range("I4").Value = 406
$ws.Range("J4").Value = 450
$ws.Range("K4").Value = 406
$ws.Range("L4").Value = 450
$ws.Range("M4").Value = -291
$ws.Range("N4").Value = -680
$ws.Range("H22").Value = 126142.875
$ws.Range("I22").Value = 143449
$ws.Range("K22").Value = 143449
$ws.Range("M22").Value = -143276
$ws.Range("H69").Value = 116250
$ws.Range("J69").Value = 136666.67
$ws.Range("L69").Value = 136666.67
$ws.Range("N69").Value = -138288.67
$ws.Range("H72").Value = 116250
$ws.Range("J72").Value = 136666.67
$ws.Range("L72").Value = 410000.01
$ws.Range("N72").Value = -418112.01
$ws.Range("H76").Value = 50166.332
$ws.Range("J76").Value = 50166.332
$ws.Range("L76").Value = 50166.332
$ws.Range("N76").Value = -50796.332
$ws.Range("H79").Value = 50166.332
$ws.Range("J79").Value = 50166.332
$ws.Range("L79").Value = 50166.332
$ws.Range("N79").Value = -52350.332
$ws.Range("H86").Value = 2308.0386
$ws.Range("I86").Value = 1364.7
$ws.Range("J86").Value = 2897.625
$ws.Range("K86").Value = 1364.7
$ws.Range("L86").Value = 2897.625
$ws.Range("M86").Value = -241.7
$ws.Range("N86").Value = -5143.625
$ws.Range("H89").Value = 2308.0386
$ws.Range("I89").Value = 1364.7
$ws.Range("J89").Value = 2897.625
$ws.Range("K89").Value = 6823.5
$ws.Range("L89").Value = 14488.125
$ws.Range("M89").Value = -1207.5
$ws.Range("N89").Value = -25720.125
$ws.Range("H125").Value = 91999.5
$ws.Range("J125").Value = 91999.5
$ws.Range("L125").Value = 91999.5
$ws.Range("N125").Value = -101839.5
$ws.Range("H130").Value = 79995
$ws.Range("J130").Value = 79995
$ws.Range("L130").Value = 79995
$ws.Range("N130").Value = -90035
$ws.Range("H134").Value = 1916.4348
$ws.Range("I134").Value = 1916.4348
$ws.Range("K134").Value = 5749.3044
$ws.Range("M134").Value = -3214.3044

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1247.7709
$ws.Range("I31").Value = 1151.25
$ws.Range("J31").Value = 1537.3334
$ws.Range("K31").Value = 1151.25
$ws.Range("L31").Value = 1537.3334
$ws.Range("M31").Value = -856.25
$ws.Range("N31").Value = -2127.3334
$ws.Range("H34").Value = 1247.7709
$ws.Range("I34").Value = 1151.25
$ws.Range("J34").Value = 1537.3334
$ws.Range("K34").Value = 1151.25
$ws.Range("L34").Value = 1537.3334
$ws.Range("M34").Value = -949.25
$ws.Range("N34").Value = -1941.3334
$ws.Range("H54").Value = 15881.6
$ws.Range("I54").Value = 14085.947
$ws.Range("J54").Value = 49999
$ws.Range("K54").Value = 14085.947
$ws.Range("L54").Value = 49999
$ws.Range("M54").Value = -13427.947
$ws.Range("N54").Value = -51315
$ws.Range("H58").Value = 37373.066
$ws.Range("I58").Value = 47682.13
$ws.Range("J58").Value = 7734.5
$ws.Range("K58").Value = 47682.13
$ws.Range("L58").Value = 7734.5
$ws.Range("M58").Value = -47479.13
$ws.Range("N58").Value = -8140.5
$ws.Range("H62").Value = 3246.5
$ws.Range("J62").Value = 3493
$ws.Range("L62").Value = 3493
$ws.Range("N62").Value = -4741
$ws.Range("H65").Value = 3246.5
$ws.Range("J65").Value = 3493
$ws.Range("L65").Value = 17465
$ws.Range("N65").Value = -23705
$ws.Range("H74").Value = 45475.7
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 45475.7
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 45475.7
$ws.Range("N74").Value = -47223.7
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 45475.7
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 45475.7
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 136427.1
$ws.Range("N77").Value = -145163.1
$ws.Range("M77").ClearContents()
$ws.Range("H94").Value = 2671.5
$ws.Range("I94").Value = 2620.6
$ws.Range("J94").Value = 2722.4
$ws.Range("K94").Value = 2620.6
$ws.Range("L94").Value = 2722.4
$ws.Range("M94").Value = -2169.6
$ws.Range("N94").Value = -3624.4
$ws.Range("H106").Value = 29950
$ws.Range("J106").Value = 29950
$ws.Range("L106").Value = 29950
$ws.Range("N106").Value = -32474
$ws.Range("H132").Value = 2088.889
$ws.Range("I132").Value = 2037.5
$ws.Range("K132").Value = 6112.5
$ws.Range("M132").Value = -3582.5
$ws.Range("H136").Value = 37373.066
$ws.Range("I136").Value = 47682.13
$ws.Range("J136").Value = 7734.5
$ws.Range("K136").Value = 143046.39
$ws.Range("L136").Value = 23203.5
$ws.Range("M136").Value = -140496.39
$ws.Range("N136").Value = -28303.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 723.43475
$ws.Range("I5").Value = 462.33334
$ws.Range("J5").Value = 1008.2727
$ws.Range("K5").Value = 1387.00002
$ws.Range("L5").Value = 3024.8181
$ws.Range("M5").Value = -1275.00002
$ws.Range("N5").Value = -3248.8181
$ws.Range("H23").Value = 439.2857
$ws.Range("J23").Value = 431.5
$ws.Range("L23").Value = 1294.5
$ws.Range("N23").Value = -1764.5
$ws.Range("H26").Value = 239.4
$ws.Range("I26").Value = 187
$ws.Range("K26").Value = 561
$ws.Range("M26").Value = -273
$ws.Range("H34").Value = 1419.6
$ws.Range("J34").Value = 1999.5
$ws.Range("L34").Value = 5998.5
$ws.Range("N34").Value = -6166.5
$ws.Range("H92").Value = 4919
$ws.Range("I92").Value = 478.66666
$ws.Range("J92").Value = 8249.25
$ws.Range("K92").Value = 1435.99998
$ws.Range("L92").Value = 24747.75
$ws.Range("M92").Value = -187.9999800000001
$ws.Range("N92").Value = -27243.75
$ws.Range("H97").Value = 507.27274
$ws.Range("I97").Value = 397.5
$ws.Range("J97").Value = 531.6667
$ws.Range("K97").Value = 1192.5
$ws.Range("L97").Value = 1595.0001
$ws.Range("M97").Value = -696.5
$ws.Range("N97").Value = -2587.0001
$ws.Range("H98").Value = 1771.6666
$ws.Range("I98").Value = 1112.3334
$ws.Range("K98").Value = 3337.0002
$ws.Range("M98").Value = -1839.0002
$ws.Range("H103").Value = 2124.25
$ws.Range("I103").Value = 197
$ws.Range("J103").Value = 2766.6667
$ws.Range("K103").Value = 591
$ws.Range("L103").Value = 8300.000100000001
$ws.Range("M103").Value = 288
$ws.Range("N103").Value = -10058.0001
$ws.Range("H107").Value = 1763.6
$ws.Range("J107").Value = 1896.2693
$ws.Range("L107").Value = 5688.8079
$ws.Range("N107").Value = -9528.8079
$ws.Range("H132").Value = 3601.4
$ws.Range("I132").Value = 3850.5
$ws.Range("J132").Value = 2605
$ws.Range("K132").Value = 34654.5
$ws.Range("L132").Value = 23445
$ws.Range("M132").Value = -32124.5
$ws.Range("N132").Value = -28505
$ws.Range("H135").Value = 723.43475
$ws.Range("I135").Value = 462.33334
$ws.Range("J135").Value = 1008.2727
$ws.Range("K135").Value = 4161.00006
$ws.Range("L135").Value = 9074.454299999999
$ws.Range("M135").Value = -1626.00006
$ws.Range("N135").Value = -14144.4543
$ws.Range("H139").Value = 3503.2
$ws.Range("I139").Value = 1114.6666
$ws.Range("K139").Value = 3343.9998
$ws.Range("M139").Value = 1796.0002
$ws.Range("H140").Value = 2118.077
$ws.Range("I140").Value = 1751.5416
$ws.Range("K140").Value = 5254.6248
$ws.Range("M140").Value = -74.6247999999996

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 29999
$ws.Range("J35").Value = 29999
$ws.Range("L35").Value = 29999
$ws.Range("N35").Value = -30595
$ws.Range("H40").Value = 35000
$ws.Range("I40").Value = 20000
$ws.Range("J40").Value = 50000
$ws.Range("K40").Value = 20000
$ws.Range("L40").Value = 50000
$ws.Range("M40").Value = -19849
$ws.Range("N40").Value = -50302
$ws.Range("H43").Value = 7500
$ws.Range("I43").Value = 7500
$ws.Range("K43").Value = 7500
$ws.Range("M43").Value = -7349
$ws.Range("H46").Value = 9760.25
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 5000
$ws.Range("N46").Value = -5312
$ws.Range("H55").Value = 8523.375
$ws.Range("I55").Value = 4031.3333
$ws.Range("K55").Value = 4031.3333
$ws.Range("M55").Value = -3704.3333
$ws.Range("H59").Value = 4452
$ws.Range("I59").Value = 4452
$ws.Range("K59").Value = 4452
$ws.Range("M59").Value = -3869
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H80").Value = 5164.8423
$ws.Range("I80").Value = 5322
$ws.Range("J80").Value = 5050.5454
$ws.Range("K80").Value = 5322
$ws.Range("L80").Value = 5050.5454
$ws.Range("M80").Value = -4324
$ws.Range("N80").Value = -7046.5454
$ws.Range("H83").Value = 5164.8423
$ws.Range("I83").Value = 5322
$ws.Range("J83").Value = 5050.5454
$ws.Range("K83").Value = 26610
$ws.Range("L83").Value = 25252.727
$ws.Range("M83").Value = -21618
$ws.Range("N83").Value = -35236.727
$ws.Range("H122").Value = 3240.1667
$ws.Range("I122").Value = 2416.0715
$ws.Range("J122").Value = 6124.5
$ws.Range("K122").Value = 7248.2145
$ws.Range("L122").Value = 18373.5
$ws.Range("M122").Value = -4798.2145
$ws.Range("N122").Value = -23273.5
$ws.Range("H132").Value = 35532.484
$ws.Range("I132").Value = 48664.13
$ws.Range("J132").Value = 5329.7
$ws.Range("K132").Value = 145992.39
$ws.Range("L132").Value = 15989.1
$ws.Range("M132").Value = -143462.39
$ws.Range("N132").Value = -21049.1

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 599.2
$ws.Range("I55").Value = 545.13635
$ws.Range("K55").Value = 545.13635
$ws.Range("M55").Value = -372.13635
$ws.Range("H132").Value = 40252.062
$ws.Range("I132").Value = 68599.89
$ws.Range("J132").Value = 3804.8572
$ws.Range("K132").Value = 205799.67
$ws.Range("L132").Value = 11414.5716
$ws.Range("M132").Value = -203269.67
$ws.Range("N132").Value = -16474.5716
$ws.Range("H136").Value = 2969.5454
$ws.Range("I136").Value = 2458.25
$ws.Range("J136").Value = 4333
$ws.Range("K136").Value = 7374.75
$ws.Range("L136").Value = 12999
$ws.Range("M136").Value = -4824.75
$ws.Range("N136").Value = -18099

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 20899.8
$ws.Range("I3").Value = 1250
$ws.Range("K3").Value = 1250
$ws.Range("M3").Value = -1136
$ws.Range("H4").Value = 19997.25
$ws.Range("J4").Value = 19999.334
$ws.Range("L4").Value = 19999.334
$ws.Range("N4").Value = -20225.334
$ws.Range("H52").Value = 25304
$ws.Range("I52").Value = 20337.889
$ws.Range("J52").Value = 69999
$ws.Range("K52").Value = 20337.889
$ws.Range("L52").Value = 69999
$ws.Range("M52").Value = -20111.889
$ws.Range("N52").Value = -70451
$ws.Range("H75").Value = 54449.8
$ws.Range("I75").Value = 50663
$ws.Range("K75").Value = 50663
$ws.Range("M75").Value = -49727
$ws.Range("H78").Value = 54449.8
$ws.Range("I78").Value = 50663
$ws.Range("K78").Value = 151989
$ws.Range("M78").Value = -147309
$ws.Range("H100").Value = 968.3333
$ws.Range("I100").Value = 817.5
$ws.Range("K100").Value = 1635
$ws.Range("M100").Value = -1094
$ws.Range("H122").Value = 5689.2
$ws.Range("I122").Value = 5487.4165
$ws.Range("K122").Value = 16462.2495
$ws.Range("M122").Value = -14012.2495
$ws.Range("H123").Value = 48000
$ws.Range("J123").Value = 48000
$ws.Range("L123").Value = 48000
$ws.Range("N123").Value = -57800
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 105421.7
$ws.Range("I126").Value = 206060.8
$ws.Range("J126").Value = 4782.6
$ws.Range("K126").Value = 618182.3999999999
$ws.Range("L126").Value = 14347.8
$ws.Range("M126").Value = -615712.3999999999
$ws.Range("N126").Value = -19287.8
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 39727.406
$ws.Range("I132").Value = 50223
$ws.Range("K132").Value = 150669
$ws.Range("M132").Value = -148139
$ws.Range("H136").Value = 3586.85
$ws.Range("I136").Value = 2639.875
$ws.Range("K136").Value = 7919.625
$ws.Range("M136").Value = -5369.625
